# Saldo_guide.xlsx update
# - Rename the worksheet to reflect the new export timestamp
# - Shift the "Dt. Referencia" (column G) date forward by one day for every data row
# - Correct the Saldo Previsto / Vl. Total values for three client rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and the corresponding workbook sheet entry) to the new export name
$ws.Name = "IClientBalance-20240614-090855-"

# Bump the reference date for every row of data (rows 2-257) from 2024-06-13 to 2024-06-14
$ws.Range("G2:G257").Value = 45457

# Fix the Saldo Previsto (D) / Vl. Total (H) figures that were corrected in this upload
$ws.Range("D101").Value = 15386.67
$ws.Range("H101").Value = 15386.67

$ws.Range("D103").Value = 15072.73
$ws.Range("H103").Value = 15072.73

$ws.Range("D107").Value = 56.27
$ws.Range("H107").Value = 56.27
